$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '38.746.14'
$ws.Range("E2").Value = '  +1.12%  '
$ws.Range("D3").Value = '2.097.76'
$ws.Range("E3").Value = '  +0.05%  '
$ws.Range("E4").Value = '  -0.09%  '
$ws.Range("D5").Value = '228.70'
$ws.Range("E5").Value = '  -0.02%  '
$ws.Range("E6").Value = '  +0.61%  '
$ws.Range("D7").Value = '62.33'
$ws.Range("E7").Value = '  +1.80%  '
$ws.Range("E8").Value = '  -0.05%  '
$ws.Range("E9").Value = '  +1.95%  '
$ws.Range("E10").Value = '  -0.59%  '
$ws.Range("E11").Value = '  -0.10%  '
$ws.Range("D12").Value = '15.80'
$ws.Range("E12").Value = '  +7.34%  '
$ws.Range("D13").Value = '2.408.69'
$ws.Range("E13").Value = '  -0.10%  '
$ws.Range("E14").Value = '  -0.82%  '
$ws.Range("E15").Value = '  +3.92%  '
$ws.Range("D16").Value = '5.51'
$ws.Range("E16").Value = '  +0.58%  '
$ws.Range("D17").Value = '2.091.99'
$ws.Range("E17").Value = '  +0.22%  '
$ws.Range("D18").Value = '38.734.86'
$ws.Range("E18").Value = '  +1.33%  '
$ws.Range("D19").Value = '71.98'
$ws.Range("E19").Value = '  +2.29%  '
$ws.Range("E20").Value = '  +0.38%  '
$ws.Range("D21").Value = '0.0₃0839'
$ws.Range("E21").Value = '  +0.75%  '
$ws.Range("D22").Value = '227.79'
$ws.Range("E22").Value = '  +1.45%  '
$ws.Range("E23").Value = '  +0.03%  '
$ws.Range("D24").Value = '2.37'
$ws.Range("E24").Value = '  -3.00%  '
$ws.Range("E25").Value = '  +0.64%  '
$ws.Range("D26").Value = '171.93'
$ws.Range("E26").Value = '  +1.24%  '
$ws.Range("E27").Value = '  +1.81%  '
$ws.Range("D28").Value = '0.139'
$ws.Range("E28").Value = '  +6.47%  '
$ws.Range("E29").Value = '  +4.29%  '
$ws.Range("E30").Value = '  +1.74%  '
$ws.Range("E31").Value = '  +3.81%  '
$ws.Range("E32").Value = '  +0.86%  '
$ws.Range("D33").Value = '4.54'
$ws.Range("E33").Value = '  +2.39%  '
$ws.Range("E34").Value = '  +0.74%  '
$ws.Range("D35").Value = '0.0619'
$ws.Range("E35").Value = '  +2.41%  '
$ws.Range("E36").Value = '  +2.82%  '
$ws.Range("D37").Value = '2.42'
$ws.Range("E37").Value = '  +1.18%  '
$ws.Range("E38").Value = '  +0.28%  '
$ws.Range("E39").Value = '  +0.05%  '
$ws.Range("E40").Value = '  +0.96%  '
$ws.Range("E41").Value = '  +4.25%  '
$ws.Range("D42").Value = '102.10'
$ws.Range("E42").Value = '  +2.25%  '
$ws.Range("D43").Value = '1.533.13'
$ws.Range("E43").Value = '  -0.83%  '
$ws.Range("E44").Value = '  -0.63%  '
$ws.Range("D45").Value = '7.80'
$ws.Range("E45").Value = '  +3.77%  '
$ws.Range("D46").Value = '0.0910'
$ws.Range("E46").Value = '  +0.19%  '
$ws.Range("E47").Value = '  +2.24%  '
$ws.Range("D48").Value = '4.13'
$ws.Range("E48").Value = '  -1.06%  '
$ws.Range("E49").Value = '  +1.43%  '
$ws.Range("E50").Value = '  -0.99%  '
$ws.Range("D51").Value = '2.291.86'
$ws.Range("E51").Value = '  -0.22%  '
